$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 (employee #3): Tran Van / A / 15-01-1990 -> Tan Thien / Lang / 15-01-2000,
# role changes from "Nhan vien quan ly kho" to "Nhan vien ban hang", final salary drops.
$ws.Range("B4").Value = "Tần Thiên"
$ws.Range("C4").Value = "Lang"
$ws.Range("D4").Value = "15-01-2000"
$ws.Range("E4").Value = "Nhân viên bán hàng"
$ws.Range("G4").Value = "6.000.000"

# Row 5: role changes to "Nhan vien ban hang", final salary updated
$ws.Range("E5").Value = "Nhân viên bán hàng"
$ws.Range("G5").Value = "7.000.000"

# Row 6: role changes to "Nhan vien ban hang", final salary updated
$ws.Range("E6").Value = "Nhân viên bán hàng"
$ws.Range("G6").Value = "8.000.000"

# Row 7: role changes to "Nhan vien ban hang", final salary updated
$ws.Range("E7").Value = "Nhân viên bán hàng"
$ws.Range("G7").Value = "9.000.000"

# Row 8: role changes to "Quan ly chuc vu", final salary updated
$ws.Range("E8").Value = "Quản lý chức vụ"
$ws.Range("G8").Value = "7.840.000"

# Row 9: role changes to "Quan ly cua hang", final salary updated
$ws.Range("E9").Value = "Quản lý cửa hàng"
$ws.Range("G9").Value = "8.320.000"

# Row 10: role changes to "Nhan vien cham soc khach hang", final salary updated
$ws.Range("E10").Value = "Nhân viên chăm sóc khách hàng"
$ws.Range("G10").Value = "6.380.000"

# Row 13: status changes from "Ngung hoat dong" to "Hoat dong"
$ws.Range("H13").Value = "Hoạt động"

# Column E (Role Name) now holds longer role-name text (e.g. "Nhan vien cham soc
# khach hang"), so it needs to grow to fit the new longest entry.
$ws.Columns.Item(5).ColumnWidth = 28.8
